$d = $word.ActiveDocument

# --- Edit 1 -----------------------------------------------------------
# The "lifecycle" paragraph is split across six runs (render()/ComponentDidMount
# snippets each in their own <w:r>). Collapse it back into a single run by
# doing a Find/Replace over the exact (already concatenated) text - Word's
# Find & Replace merges every run the match spans into one run holding the
# replacement text.
$lifecycleText = "Trong quá trình lifecycle , khi một file .js được gọi , nó sẽ chạy hàm render(){} trước sau đó mới tới ComponentDidMount(){} , sau đó nó sẽ chạy luôn phiên render(){} => ComponentDidMount(){}"
$d.Content.Find.Execute($lifecycleText, $true, $false, $false, $false, $false, $true, 1, $false, $lifecycleText, 2) | Out-Null

# --- Edit 2 -----------------------------------------------------------
# Remove the three empty paragraphs that sit between the "&nbsp hoặc <></>"
# code line and the "npm toastify" Heading2. Locate them by walking the
# paragraph collection instead of hard-coded indices, in case earlier edits
# shifted paragraph numbering.
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*nbsp*") {
        $first = $d.Paragraphs.Item($i + 1)
        $last = $d.Paragraphs.Item($i + 3)
        $rng = $d.Range($first.Range.Start, $last.Range.End)
        $rng.Delete()
        break
    }
}
